# Change TestObject folder structure and page naming convention to separate
# from functional tests: prefix each page's sheet name with "V_" and move the
# active tab/selection from the ResourcesPage to the AboutAboutPage.

$wb = $excel.ActiveWorkbook

# Rename the sheets to the new "V_" naming convention.
$wb.Worksheets.Item("HomePage").Name = "V_HomePage"
$wb.Worksheets.Item("DataModelPage").Name = "V_DataModelPage"
$wb.Worksheets.Item("ResourcesPage").Name = "V_ResourcesPage"
$wb.Worksheets.Item("AboutAboutPage").Name = "V_AboutAboutPage"

# Move the active/selected tab from the Resources page to the About page,
# and update the selected cell on the About page.
$wsAbout = $wb.Worksheets.Item("V_AboutAboutPage")
$wsAbout.Activate()
$wsAbout.Range("B33").Select()
